$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename column C header: "Udział dnia" -> "Udział dnia w godzinach" ---
$ws.Range("C2").Value = "Udział dnia w godzinach"

# --- Add new column E: "Udział dnia w TO" (share of the day in the monthly Turn Over) ---

# Clone the header formatting (style) from D2 and set the new header text.
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value = "Udział dnia w TO"

# Clone the data-cell formatting (percentage style) from column C down column E.
$ws.Range("C4:C34").Copy()
$ws.Range("E4:E34").PasteSpecial(-4122)

# Clone the totals-row formatting from B35 for the E35 total cell.
$ws.Range("B35").Copy()
$ws.Range("E35").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Calculate and write each day's share of the monthly Turn Over: D{row} / D35 (the
# monthly Turn Over total), mirroring the new "daily share in monthly turn over" logic.
$monthlyTurnOver = $ws.Cells.Item(35, 4).Value()

for ($r = 4; $r -le 34; $r++) {
    $dailyTurnOver = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 5).Value = $dailyTurnOver / $monthlyTurnOver
}

# Total row: sum of the daily shares equals 100% of the monthly Turn Over.
$ws.Range("E35").Value = 1.0

# --- Column widths: widen C (longer header) and size the new E column ---
$ws.Columns.Item(3).ColumnWidth = 23.0625
$ws.Columns.Item(5).ColumnWidth = 16.44921875
